$d = $word.ActiveDocument

# 1. Find the run that currently ends the Title paragraph ("...Co. ") and
#    replace its trailing non-breaking space with ": " so it can be joined
#    with the subtitle text that follows.
$d.Content.Find.Execute(
    "Let’s talk about Thurstone & Co. ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Let’s talk about Thurstone & Co.: An information-theoretical model for comparative judgments, and its statistical translation",
    2
) | Out-Null

# 2. Remove the now-orphaned Subtitle paragraph (style "Subtitle"), including
#    its paragraph mark, leaving only the merged Title paragraph.
foreach ($p in @($d.Paragraphs)) {
    if ($p.Style.NameLocal -eq "Subtitle") {
        $p.Range.Delete()
    }
}
